# Insert a new data row at row 51 (shifting existing rows 51-74 down to 52-75)
# for the weekly Fruta/Hortalizas price update ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51").Insert()

$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44489
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100102
$ws.Range("H51").Value = "Cítricos"
$ws.Range("I51").Value = 100102004
$ws.Range("J51").Value = "Mandarina"
$ws.Range("K51").Value = "Murcott"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 250
$ws.Range("N51").Value = 14000
$ws.Range("O51").Value = 15000
$ws.Range("P51").Value = 14500
$ws.Range("Q51").Value = "$/caja 20 kilos"
$ws.Range("R51").Value = "Región de Coquimbo"
$ws.Range("S51").Value = 725
$ws.Range("T51").Value = 20
